# Insert a new L1 seed row ("L1_MinimumBiasHF0", index 462) into the
# prescale table directly above the existing "L1_NotBptxOR" row (357),
# pushing it and all rows below it down by one (357 -> 358, ..., 391 -> 392).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 357 and everything below it down by one row, leaving a blank
# row 357 ready to be populated.
$ws.Rows(357).Insert()

# Populate the newly-inserted row with the new prefiring monitoring seed.
$ws.Range("A357").Value = 462
$ws.Range("B357").Value = "L1_MinimumBiasHF0"
$ws.Range("C357:K357").Value = 0
